# edit.ps1 - applies the DDO-OFDM_I_Tx_pt-br.pptx slide edits:
#   1) "TextBox 1606": merge the two runs "Amostras no " + "Domínio da
#      Frequência" into a single run "Amostras no Domínio da Frequência"
#      (text is unchanged once concatenated, but PowerPoint's COM layer
#      only rewrites the run structure when the assigned text actually
#      differs from the current text, so we first stomp it with a throw-
#      away value). The shape is a spAutoFit, vertical text box, so
#      rewriting its text can nudge its autofit Height; we restore the
#      original EMU-exact size afterwards.
#   2) "TextBox 2321": change the caption text "Vetor Símbolo OFDM" to
#      "Vetor de Símbolos OFDM" and reposition/resize it from
#      off(8208131,868789)/ext(1797608,313484) to
#      off(8064115,845952)/ext(2147063,313484) (EMU).
#
# Shape.Left/Top/Width/Height are expressed in points and are backed by
# 32-bit floats, so naively writing target_emu/12700.0 can truncate to the
# neighbouring EMU once re-saved. EmuToPts searches for a point value that
# round-trips (via the same float32 truncation PowerPoint performs) to the
# exact target EMU.

function EmuToPts($targetEmu) {
    $basePts = [double]$targetEmu / 12700.0
    for ($i = -2000; $i -le 2000; $i++) {
        $cand = $basePts + ($i * 0.0000001)
        $f = [double]([float]$cand)
        $e = [int64]($f * 12700.0)
        if ($e -eq $targetEmu) {
            return $cand
        }
    }
    return $basePts
}

function Get-ShapeByName($shapes, $name) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) "Amostras no " + "Domínio da Frequência" -> single run ---
$tbFreq = Get-ShapeByName $s.Shapes "TextBox 1606"
$origHeight = $tbFreq.Height
# force a genuine text change so the run split collapses into one run
$tbFreq.TextFrame.TextRange.Text = "~"
$tbFreq.TextFrame.TextRange.Text = "Amostras no Domínio da Frequência"
# restore the exact original autosized height (EMU-accurate)
$tbFreq.Height = $origHeight

# --- 2) "Vetor Símbolo OFDM" -> "Vetor de Símbolos OFDM" + reposition ---
$tbVetor = Get-ShapeByName $s.Shapes "TextBox 2321"
$tbVetor.TextFrame.TextRange.Text = "Vetor de Símbolos OFDM"
$tbVetor.Left = EmuToPts 8064115
$tbVetor.Top = EmuToPts 845952
$tbVetor.Width = EmuToPts 2147063
$tbVetor.Height = EmuToPts 313484
